$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "247.05"
Set-TextCell 2 7 "16"
Set-TextCell 3 4 "21.64"
Set-TextCell 3 7 "16"
Set-TextCell 4 4 "5.417"
Set-TextCell 4 7 "16"
Set-TextCell 5 4 "0.05682"
Set-TextCell 5 7 "16"
Set-TextCell 6 4 "3.381"
Set-TextCell 6 7 "16"
Set-TextCell 7 4 "0.8093"
Set-TextCell 7 7 "16"
Set-TextCell 8 4 "1.019"
Set-TextCell 8 7 "16"
Set-TextCell 9 4 "0.1448"
Set-TextCell 9 7 "16"
Set-TextCell 10 4 "0.07510"
Set-TextCell 10 7 "16"
Set-TextCell 11 4 "0.03145"
Set-TextCell 11 7 "16"
Set-TextCell 12 2 "BitrueCoin"
Set-TextCell 12 3 "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell 12 4 "0.03043"
Set-TextCell 12 5 "11BitrueCoinBTR"
Set-TextCell 12 7 "16"
Set-TextCell 13 2 "BitMartToken"
Set-TextCell 13 3 "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell 13 4 "0.09262"
Set-TextCell 13 5 "12BitMartTokenBMX"
Set-TextCell 13 7 "16"
Set-TextCell 14 2 "MCDex"
Set-TextCell 14 3 "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell 14 4 "3.610"
Set-TextCell 14 5 "13MCDexMCBBestin24h"
Set-TextCell 14 7 "16"
Set-TextCell 15 2 "BitForexToken"
Set-TextCell 15 3 "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell 15 4 "0.001657"
Set-TextCell 15 5 "14BitForexTokenBF"
Set-TextCell 15 7 "16"
Set-TextCell 16 2 "CoinExToken"
Set-TextCell 16 3 "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell 16 4 "0.04710"
Set-TextCell 16 5 "15CoinExTokenCET"
Set-TextCell 16 7 "16"
Set-TextCell 17 2 "One"
Set-TextCell 17 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell 17 4 "0.0005866"
Set-TextCell 17 5 "16OneONE"
Set-TextCell 17 7 "16"
Set-TextCell 18 2 "TigerCash"
Set-TextCell 18 3 "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell 18 4 "0.006349"
Set-TextCell 18 5 "17TigerCashTCH"
Set-TextCell 18 7 "16"
Set-TextCell 19 2 "HotbitToken"
Set-TextCell 19 3 "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell 19 4 "0.005010"
Set-TextCell 19 5 "18HotbitTokenHTB"
Set-TextCell 19 7 "16"
Set-TextCell 20 2 "BitKan"
Set-TextCell 20 3 "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell 20 4 "0.001043"
Set-TextCell 20 5 "19BitKanKAN"
Set-TextCell 20 7 "16"
Set-TextCell 21 2 "NitroEx"
Set-TextCell 21 3 "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell 21 4 "0.0001501"
Set-TextCell 21 5 "20NitroExNTX"
Set-TextCell 21 7 "16"
Set-TextCell 22 2 "UpBots"
Set-TextCell 22 3 "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextCell 22 4 "0.0003103"
Set-TextCell 22 5 "21UpBotsUBXT"
Set-TextCell 22 7 "16"
Set-TextCell 23 2 "LEO"
Set-TextCell 23 3 "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell 23 4 "3.768"
Set-TextCell 23 5 "22LEOLEO"
Set-TextCell 23 7 "16"
Set-TextCell 24 2 "KuCoinToken"
Set-TextCell 24 3 "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell 24 4 "6.409"
Set-TextCell 24 5 "23KuCoinTokenKCS"
Set-TextCell 24 7 "16"
Set-TextCell 25 2 "BTSEToken"
Set-TextCell 25 3 "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell 25 4 "2.103"
Set-TextCell 25 5 "24BTSETokenBTSE"
Set-TextCell 25 7 "16"
Set-TextCell 26 2 "BitpandaEcosystemToken"
Set-TextCell 26 3 "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell 26 4 "0.3287"
Set-TextCell 26 5 "25BitpandaEcosystemTokenBEST"
Set-TextCell 26 7 "16"
Set-TextCell 27 2 "ProBitToken"
Set-TextCell 27 3 "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell 27 4 "0.1305"
Set-TextCell 27 5 "26ProBitTokenPROB"
Set-TextCell 27 7 "16"
Set-TextCell 28 7 "16"
Set-TextCell 29 7 "16"
Set-TextCell 30 7 "16"
Set-TextCell 31 7 "16"
Set-TextCell 32 7 "16"
Set-TextCell 33 7 "16"
Set-TextCell 34 7 "16"
Set-TextCell 35 7 "16"
Set-TextCell 36 7 "16"
Set-TextCell 37 7 "16"
Set-TextCell 38 7 "16"
Set-TextCell 39 7 "16"
Set-TextCell 40 4 "0.04038"
Set-TextCell 40 7 "16"
Set-TextCell 41 4 "0.006967"
Set-TextCell 41 7 "16"
Set-TextCell 42 7 "16"
Set-TextCell 43 4 "0.002932"
Set-TextCell 43 7 "16"
Set-TextCell 44 4 "0.008506"
Set-TextCell 44 7 "16"
Set-TextCell 45 4 "0.00005939"
Set-TextCell 45 7 "16"
Set-TextCell 46 7 "16"
Set-TextCell 47 4 "0.0005505"
Set-TextCell 47 7 "16"
Set-TextCell 48 4 "0.6831"
Set-TextCell 48 7 "16"
Set-TextCell 49 4 "0.01459"
Set-TextCell 49 7 "16"
Set-TextCell 50 4 "0.00002102"
Set-TextCell 50 7 "16"
Set-TextCell 51 7 "16"
